$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp text update ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 14:35"

# --- Rows 26-28: Suecia overtakes Catar and Portugal in the ranking ---
$ws.Range("A26").Value = "Suecia"
$ws.Range("B26").Value = 28582
$ws.Range("C26").Value = 673
$ws.Range("D26").Value = 4971
$ws.Range("E26").Value = 20082
$ws.Range("F26").Value = 351
$ws.Range("G26").Value = 69
$ws.Range("H26").Value = 3529

$ws.Range("A27").Value = "Catar"
$ws.Range("B27").Value = 28272
$ws.Range("C27").Value = 1733
$ws.Range("D27").Value = 3356
$ws.Range("E27").Value = 24902
$ws.Range("F27").Value = 72
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 14

$ws.Range("A28").Value = "Portugal"
$ws.Range("B28").Value = 28132
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 3182
$ws.Range("E28").Value = 23775
$ws.Range("F28").Value = 103
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 1175

# --- Row 47: Dinamarca data refresh ---
$ws.Range("D47").Value = 8805
$ws.Range("E47").Value = 1371
$ws.Range("F47").Value = 35
$ws.Range("G47").Value = 4
$ws.Range("H47").Value = 537

# --- Row 65: data refresh ---
$ws.Range("F65").Value = 31

# --- Row 75: data refresh ---
$ws.Range("B75").Value = 2636
$ws.Range("C75").Value = 24
$ws.Range("D75").Value = 2136
$ws.Range("E75").Value = 489

# --- Row 78: data refresh ---
$ws.Range("B78").Value = 2221
$ws.Range("C78").Value = 8
$ws.Range("D78").Value = 1850
$ws.Range("E78").Value = 277
$ws.Range("F78").Value = 7

# --- Rows 94-95: Consejo Danes para los Refugiados overtakes Somalia ---
$ws.Range("A94").Value = "Consejo Danes para los Refugiados"
$ws.Range("B94").Value = 1242
$ws.Range("C94").Value = 73
$ws.Range("D94").Value = 157
$ws.Range("E94").Value = 1035
$ws.Range("F94").Value = 0
$ws.Range("H94").Value = 50

$ws.Range("A95").Value = "Somalia"
$ws.Range("B95").Value = 1219
$ws.Range("D95").Value = 130
$ws.Range("E95").Value = 1037
$ws.Range("F95").Value = 2
$ws.Range("H95").Value = 52

# --- Rows 193-195: Santa Lucia & Belice overtake Nueva Caledonia ---
$ws.Range("A193").Value = "Santa Lucia"

$ws.Range("A194").Value = "Belice"
$ws.Range("D194").Value = 16
$ws.Range("H194").Value = 2

$ws.Range("A195").Value = "Nueva Caledonia"
$ws.Range("D195").Value = 18
$ws.Range("H195").Value = 0

# --- Rows 200-201: Curazao overtakes Dominica ---
$ws.Range("A200").Value = "Curazao"
$ws.Range("D200").Value = 14
$ws.Range("H200").Value = 1

$ws.Range("A201").Value = "Dominica"
$ws.Range("D201").Value = 15
$ws.Range("H201").Value = 0
